$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.747119
$ws.Range("H2").Value = 2.241357
$ws.Range("I2").Value = 0.03096954854571248
$ws.Range("J2").Value = 0.03096954854571248
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 297.8183156666666
$ws.Range("N2").Value = 893.4549469999999
$ws.Range("O2").Value = 0.8852156413092672
$ws.Range("P2").Value = 0.8852156413092673
$ws.Range("Q2").Value = 222.5057221825643
$ws.Range("R2").Value = 2002.551499643079
$ws.Range("S2").Value = 0.02741472877695136
$ws.Range("T2").Value = 0.02741472877695136
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.747119
$ws.Range("H3").Value = 2.241357
$ws.Range("I3").Value = 0.03096954854571248
$ws.Range("J3").Value = 0.03096954854571248
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 24.34034433333333
$ws.Range("N3").Value = 73.021033
$ws.Range("O3").Value = 0.07234764413494278
$ws.Range("P3").Value = 0.0723476441349428
$ws.Range("Q3").Value = 18.18513371797567
$ws.Range("R3").Value = 163.666203461781
$ws.Range("S3").Value = 0.002240573877205041
$ws.Range("T3").Value = 0.002240573877205042
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.747119
$ws.Range("H4").Value = 2.241357
$ws.Range("I4").Value = 0.03096954854571248
$ws.Range("J4").Value = 0.03096954854571248
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.277234
$ws.Range("N4").Value = 42.831702
$ws.Range("O4").Value = 0.04243671455578994
$ws.Range("P4").Value = 0.04243671455578994
$ws.Range("Q4").Value = 10.666792788846
$ws.Range("R4").Value = 96.00113509961399
$ws.Range("S4").Value = 0.00131424589155608
$ws.Range("T4").Value = 0.00131424589155608
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.74619233333334
$ws.Range("H5").Value = 59.23857700000001
$ws.Range("I5").Value = 0.8185184181638298
$ws.Range("J5").Value = 0.8185184181638298
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 297.8183156666666
$ws.Range("N5").Value = 893.4549469999999
$ws.Range("O5").Value = 0.8852156413092672
$ws.Range("P5").Value = 0.8852156413092673
$ws.Range("Q5").Value = 5880.77774154338
$ws.Range("R5").Value = 52926.99967389042
$ws.Range("S5").Value = 0.7245653064583415
$ws.Range("T5").Value = 0.7245653064583416
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.74619233333334
$ws.Range("H6").Value = 59.23857700000001
$ws.Range("I6").Value = 0.8185184181638298
$ws.Range("J6").Value = 0.8185184181638298
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 24.34034433333333
$ws.Range("N6").Value = 73.021033
$ws.Range("O6").Value = 0.07234764413494278
$ws.Range("P6").Value = 0.0723476441349428
$ws.Range("Q6").Value = 480.6291206655602
$ws.Range("R6").Value = 4325.662085990041
$ws.Range("S6").Value = 0.05921787923521304
$ws.Range("T6").Value = 0.05921787923521305
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.74619233333334
$ws.Range("H7").Value = 59.23857700000001
$ws.Range("I7").Value = 0.8185184181638298
$ws.Range("J7").Value = 0.8185184181638298
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.277234
$ws.Range("N7").Value = 42.831702
$ws.Range("O7").Value = 0.04243671455578994
$ws.Range("P7").Value = 0.04243671455578994
$ws.Range("Q7").Value = 281.9210085520061
$ws.Range("R7").Value = 2537.289076968054
$ws.Range("S7").Value = 0.03473523247027515
$ws.Range("T7").Value = 0.03473523247027515
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.630999
$ws.Range("H8").Value = 10.892997
$ws.Range("I8").Value = 0.1505120332904577
$ws.Range("J8").Value = 0.1505120332904577
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 297.8183156666666
$ws.Range("N8").Value = 893.4549469999999
$ws.Range("O8").Value = 0.8852156413092672
$ws.Range("P8").Value = 0.8852156413092673
$ws.Range("Q8").Value = 1081.378006367351
$ws.Range("R8").Value = 9732.402057306157
$ws.Range("S8").Value = 0.1332356060739743
$ws.Range("T8").Value = 0.1332356060739743
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.630999
$ws.Range("H9").Value = 10.892997
$ws.Range("I9").Value = 0.1505120332904577
$ws.Range("J9").Value = 0.1505120332904577
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 24.34034433333333
$ws.Range("N9").Value = 73.021033
$ws.Range("O9").Value = 0.07234764413494278
$ws.Range("P9").Value = 0.0723476441349428
$ws.Range("Q9").Value = 88.37976593398899
$ws.Range("R9").Value = 795.4178934059009
$ws.Range("S9").Value = 0.01088919102252469
$ws.Range("T9").Value = 0.0108891910225247
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.630999
$ws.Range("H10").Value = 10.892997
$ws.Range("I10").Value = 0.1505120332904577
$ws.Range("J10").Value = 0.1505120332904577
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.277234
$ws.Range("N10").Value = 42.831702
$ws.Range("O10").Value = 0.04243671455578994
$ws.Range("P10").Value = 0.04243671455578994
$ws.Range("Q10").Value = 51.84062237676599
$ws.Range("R10").Value = 466.565601390894
$ws.Range("S10").Value = 0.006387236193958704
$ws.Range("T10").Value = 0.006387236193958707
